$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume columns to stay text so Excel does not
# auto-convert numeric-looking strings (e.g. "243.29") to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.430.31"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.867.74"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "243.29"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "0.7059"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("D10").Value = "24.48"
$ws.Range("E10").Value = "  -2.08%  "
$ws.Range("D11").Value = "0.07993"
$ws.Range("E11").Value = "  -3.80%  "
$ws.Range("D12").Value = "1.874.52"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "5.204"
$ws.Range("D14").Value = "93.29"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "0.6992"
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").Value = "6.505"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.000008378"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "29.484.22"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").Value = "252.45"
$ws.Range("E19").Value = "  +3.66%  "
$ws.Range("D20").Value = "2.137.21"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("D23").Value = "7.624"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "0.1556"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").Value = "9.007"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").Value = "161.12"
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("D28").Value = "18.72"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").Value = "4.320"
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("D32").Value = "1.204"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("D34").Value = "1.888"
$ws.Range("D35").Value = "0.7475"
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").Value = "2.715"
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("D38").Value = "0.01881"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("D39").Value = "1.271.93"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D40").Value = "2.747"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("D41").Value = "0.8941"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("D42").Value = "6.099"
$ws.Range("E42").Value = "  -5.97%  "
$ws.Range("D43").Value = "108.98"
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("D44").Value = "71.35"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "0.00000000129"
$ws.Range("E46").Value = "  -2.74%  "
$ws.Range("D47").Value = "2.037.59"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").Value = "9.569"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").Value = "1.790"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").Value = "0.5177"
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("E51").Value = "  -1.48%  "

# Restore default (unformatted) style now that values are locked in as text.
$ws.Range("D2:E51").ClearFormats()
